$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.244.32"
$ws.Range("E2").Value = "  +3.17%  "

$ws.Range("D3").Value = "1.896.78"

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Formula = "=""325.59"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  +3.55%  "

$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").Formula = "=""0.5164"""
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").Formula = "=""0.4001"""
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("D9").Formula = "=""0.08434"""
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").Formula = "=""42.62"""
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").Formula = "=""1.117"""
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Formula = "=""23.46"""
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  +13.51%  "

$ws.Range("D13").Formula = "=""6.427"""
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("E13").Value = "  +2.49%  "

$ws.Range("D14").Value = "1.894.49"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").Formula = "=""7.336"""
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Formula = "=""94.42"""
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").Formula = "=""0.06647"""
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  -1.22%  "

$ws.Range("D20").Formula = "=""18.22"""
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163) | Out-Null
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").Formula = "=""5.949"""
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  -0.99%  "

$ws.Range("D23").Value = "30.249.38"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").Formula = "=""2.225"""
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").Value = "2.108.17"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Formula = "=""21.66"""
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  +3.82%  "

$ws.Range("D28").Formula = "=""161.16"""
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = "  +1.19%  "

$ws.Range("D29").Formula = "=""2.360"""
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  -2.94%  "

$ws.Range("D30").Formula = "=""128.78"""
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +1.39%  "

$ws.Range("D31").Formula = "=""1.098"""
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  +3.66%  "

$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").Formula = "=""6.072"""
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Formula = "=""3.754"""
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("D35").Formula = "=""0.02499"""
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("D36").Formula = "=""0.06561"""
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").Formula = "=""5.276"""
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").Formula = "=""0.2202"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").Formula = "=""1.218"""
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("E40").Value = "  +4.63%  "

$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").Formula = "=""8.725"""
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("D43").Formula = "=""1.233"""
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Formula = "=""0.6104"""
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").Formula = "=""13.21"""
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("D46").Formula = "=""3.704"""
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("D47").Formula = "=""2.057"""
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").Formula = "=""1.235"""
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").Formula = "=""124.50"""
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("D51").Formula = "=""79.05"""
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  +1.77%  "
